$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate and delete the "TextBox 67" shape (text "Updates in location",
# positioned at x=345480, y=3584037) that was removed from the slide.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "TextBox 67") {
        $shape.Delete()
    }
}
